$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.488.93"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "1.879.52"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'246.03"
$ws.Range("E5").Value = "  +5.26%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "'0.4762"
$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D8").Value = "'0.2902"
$ws.Range("E8").Value = "  +1.66%  "

$ws.Range("D9").Value = "'0.06530"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("D10").Value = "'22.00"
$ws.Range("E10").Value = "  +3.74%  "

$ws.Range("D11").Value = "'0.07740"
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").Value = "'0.7434"
$ws.Range("E12").Value = "  +8.98%  "

$ws.Range("E13").Value = "  +3.55%  "

$ws.Range("D14").Value = "1.876.77"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").Value = "'5.133"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").Value = "'272.93"
$ws.Range("E16").Value = "  +1.65%  "

$ws.Range("D17").Value = "30.481.86"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").Value = "'13.66"
$ws.Range("E18").Value = "  +2.63%  "

$ws.Range("D19").Value = "'0.000007592"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "2.129.91"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "'5.239"
$ws.Range("E23").Value = "  +1.99%  "

$ws.Range("D24").Value = "'6.185"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("D25").Value = "'9.307"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "'164.46"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").Value = "'18.94"
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("D28").Value = "'1.950"
$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("E29").Value = "  +0.80%  "

$ws.Range("D30").Value = "'0.09999"
$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("D31").Value = "'1.509"
$ws.Range("E31").Value = "  +4.11%  "

$ws.Range("D32").Value = "'4.327"
$ws.Range("E32").Value = "  +2.38%  "

$ws.Range("D33").Value = "'4.065"
$ws.Range("E33").Value = "  +1.49%  "

$ws.Range("D34").Value = "'0.04786"
$ws.Range("E34").Value = "  +2.40%  "

$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("D36").Value = "'0.7002"
$ws.Range("E36").Value = "  +1.83%  "

$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").Value = "'0.01871"
$ws.Range("E38").Value = "  +2.30%  "

$ws.Range("D39").Value = "'2.741"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'6.358"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("D41").Value = "'1.937"
$ws.Range("E41").Value = "  +2.45%  "

$ws.Range("D42").Value = "'70.07"
$ws.Range("E42").Value = "  -1.25%  "

$ws.Range("D43").Value = "'0.4174"
$ws.Range("E43").Value = "  +3.04%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8402"
$ws.Range("E44").Value = "  +0.80%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "'102.92"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("D47").Value = "'9.327"
$ws.Range("E47").Value = "  +3.02%  "

$ws.Range("D48").Value = "'7.097"
$ws.Range("E48").Value = "  +2.16%  "

$ws.Range("D49").Value = "'35.45"
$ws.Range("E49").Value = "  +4.45%  "

$ws.Range("D50").Value = "'919.73"
$ws.Range("E50").Value = "  -1.59%  "

$ws.Range("D51").Value = "'0.05627"
$ws.Range("E51").Value = "  +0.97%  "
